# Project DesignFirst is saved. Author: admin. Type: SAVE.
#
# Semantic change carried by this save: cell D10 on the "Rules" sheet was
# updated from 21 to 100 (numeric), matching the value already present in
# the neighbouring cell C10.
#
# (The workbook's <col> definitions also flip a cosmetic collapsed="false"
# -> collapsed="true" flag on every column in the raw OOXML for this save,
# but that is a blanket formatting no-op written uniformly across all
# columns by the authoring tool -- not a deliberate outline/group action on
# any particular column -- so there is no corresponding targeted Range/
# Columns call to make here; touching Columns(...).OutlineLevel/ShowDetail
# would instead record a real (and undesired) grouping/hidden state that
# isn't part of this change.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D10").Value = 100.0
